$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was a placeholder "1", now the proper municipality name)
$ws.Name = "ჩხოროწყუ"

# Drop the "(census results)" note row entirely - the row below (a blank
# spacer row) shifts up and becomes the new row 2
$ws.Rows(2).Delete()

# Drop the obsolete 1989/2002 year columns; the 2014 column shifts left
# into column B
$ws.Range("B1:C1").EntireColumn.Delete()
